$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "metrics.tsv" data-dictionary rows (rows 2-8) to reflect the
# new Global Architecture for Health Emergency Preparedness, Response and
# Resilience (HEPR) metrics content, and to reorder field descriptions.

# Row 2 (metric_id)
$ws.Range("B2").Value = "metric_id"
$ws.Range("C2").Value = "A unique ID associated with the specified metric"

# Row 3 (framework)
$ws.Range("C3").Value = "The name of the framework or system in which the metric is defined, including information on the edition (e.g., ""JEE 3.0"", ""SPAR 2.0"", ""Health Emergency Preparedness, Response and Resilience (HEPR)"")"

# Row 4 (pillar)
$ws.Range("C4").Value = "The pillar or key system of global health security that the indicator/attribute correspond to; one of: Prevent, Detect, Respond, or IHR Related Hazards and Points of Entry and Border Health (for JEE SPAR), or Collaborative surveillance, Access to countermeasures, Emergency coordination, Clinical care, or Community protection (Global Architecture for Health Emergency Preparedness, Response and Resilience)"
$ws.Range("D4").Value = "For more details, please see https://www.who.int/publications/i/item/9789240051980; note that the second edition SPAR (SPAR 2.0) does not explicitly define these pillars, for SPAR metrics, pillars were inferred by the research team based on pillar-indicator pairs as defined in the JEE"

# Row 5 (capacity)
$ws.Range("C5").Value = "The capacity associated with the metric"
$ws.Range("D5").Value = "For more details, please see  JEE or SPAR reference documents"

# Row 6 (indicator)
$ws.Range("C6").Value = "The indicator associated with the metric"
$ws.Range("D6").Value = "For more details, please see  JEE or SPAR reference documents"

# Row 7 (score)
$ws.Range("C7").Value = "A numeric score that assesses country performance against the metric"
$ws.Range("D7").Value = "For more details, please see  JEE or SPAR reference documents; note that the metrics within the Health Emergency Preparedness, Response and Resilience (HEPR) are not specifically scored, so no score will be indicated, only a written description of the desired capacity "

# Row 8 (attribute)
$ws.Range("C8").Value = "The attribute that is required to obtain the specified score on the metric"
$ws.Range("D8").Value = "For more details, please see  JEE or SPAR reference documents; note that the metrics within the Health Emergency Preparedness, Response and Resilience (HEPR) are not specifically scored, so no score will be indicated, only a written description of the desired capacity "

# --- Adjust column widths (column B narrower, column C wider) to fit the
# updated text content.
$ws.Columns.Item(2).ColumnWidth = 16.33
$ws.Columns.Item(3).ColumnWidth = 73.5
